$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = 5
    4  = -4
    5  = 0
    8  = 0
    9  = 2
    10 = 6
    12 = -12
    13 = -3
    14 = -5
    16 = 0
    17 = -1
    18 = -8
    20 = -3
    21 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
